### Lesson 10 homework — add column C data, extend the chart with a 2nd series + a
### polynomial trendline on the 1st series, drop the fixed X-axis max, and move/resize
### the chart on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data in column C (the homework numbers for lesson 10)
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = 4
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 82
$ws.Range("C4").Value = 372
$ws.Range("C5").Value = 1478
$ws.Range("C6").Value = 6887
$ws.Range("C7").Value = 26633
$ws.Range("C8").Value = 96630
$ws.Range("C9").Value = 384116
$ws.Range("C10").Value = 1502695

$ws.Range("C11").Select()

# ---------------------------------------------------------------------------
# 2. Chart edits
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# 2a. Polynomial (order 2) trendline on the first series
$s1 = $chart.SeriesCollection(1)
$tl = $s1.Trendlines().Add()
$tl.Type = 3
$tl.Order = 2
$tl.DisplayEquation = $false
$tl.DisplayRSquared = $false
$tl.Border.Weight = 3.25
$tl.Border.LineStyle = 2
$tl.Border.Color = 0

# 2b. Second series plotting the new column C against column A
$sc = $chart.SeriesCollection()
$ns = $sc.NewSeries()
$ns.XValues = "=Лист1!`$A`$1:`$A`$10"
$ns.Values = "=Лист1!`$C`$1:`$C`$10"
$ns.Border.Weight = 1.5
$ns.Border.LineStyle = 1
$ns.Border.Color = 3243501

# 2c. Drop the fixed maximum on the X (value) axis so it autoscales again
$xAxis = $chart.Axes(1)
$xAxis.MaximumScaleIsAuto = $true

# 2d. Reposition / resize the chart on the sheet
$co.Left = 195.5625
$co.Top = 38.25
$co.Width = 699.75
$co.Height = 382.5
